$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Drop the "is_halal" column (N) entirely - header + all 20 data rows.
# ---------------------------------------------------------------------------
$ws.Range("N1").EntireColumn.Delete()

# ---------------------------------------------------------------------------
# 2. Replace the first two product rows (previously "Lager Beer A"/"Lager Beer B")
#    with the new products "Baileys" and "Aperol".
# ---------------------------------------------------------------------------
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = "Baileys"
$ws.Range("F4").Value = "Liquer"
$ws.Range("G4").Value = 17
$ws.Range("H4").Value = 700
$ws.Range("I4").Value = "Krim susu, whiskey Irlandia, kakao"
$ws.Range("J4").Value = 20
$ws.Range("K4").Value = 450000
$ws.Range("L4").Value = "Ireland "
$ws.Range("M4").Value = "Bottle"

$ws.Range("D5").Value = 2
$ws.Range("E5").Value = "Aperol"
$ws.Range("F5").Value = "Liquer"
$ws.Range("G5").Value = 11
$ws.Range("H5").Value = 750
$ws.Range("I5").Value = "Jeruk pahit, herbal, rhubarb"
$ws.Range("J5").Value = 16
$ws.Range("K5").Value = 350000
$ws.Range("L5").Value = "Italy"
$ws.Range("M5").Value = "Bottle"

# ---------------------------------------------------------------------------
# 3. Apply an Indonesian-Rupiah accounting number format to the whole
#    price_idr column (K4:K23).
# ---------------------------------------------------------------------------
$rupiah = '_-[$Rp-421]* #,##0_-;\-[$Rp-421]* #,##0_-;_-[$Rp-421]* "-"??_-;_-@_-'
$ws.Range("K4:K23").NumberFormat = $rupiah

# ---------------------------------------------------------------------------
# 4. Move the active selection to I8 (matches the author's saved cursor spot).
# ---------------------------------------------------------------------------
$null = $ws.Range("I8").Select()

Write-Host "done"
